# Refresh crypto price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.Value = "'42.701.90"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$dCell = $ws.Range("D3")
$dCell.Value = "'2.533.67"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.04%  "
$dCell = $ws.Range("D5")
$dCell.Value = "'315.68"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$dCell = $ws.Range("D6")
$dCell.Value = "'95.53"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  -0.05%  "
$dCell = $ws.Range("D9")
$dCell.Value = "'0.531"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "
$dCell = $ws.Range("D10")
$dCell.Value = "'35.66"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -1.75%  "
$dCell = $ws.Range("D11")
$dCell.Value = "'0.0804"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "
$dCell = $ws.Range("D12")
$dCell.Value = "'7.51"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("E13").Value = "  -2.55%  "
$dCell = $ws.Range("D14")
$dCell.Value = "'2.919.87"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "
$dCell = $ws.Range("D15")
$dCell.Value = "'2.566.34"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +2.02%  "
$dCell = $ws.Range("D16")
$dCell.Value = "'15.04"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -3.36%  "
$dCell = $ws.Range("D17")
$dCell.Value = "'0.847"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "
$dCell = $ws.Range("D18")
$dCell.Value = "'42.772.51"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +2.96%  "
$dCell = $ws.Range("D20")
$dCell.Value = "'12.76"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "
$dCell = $ws.Range("D21")
$dCell.Value = "'0.0₃0960"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$dCell = $ws.Range("D22")
$dCell.Value = "'69.62"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "
$dCell = $ws.Range("D23")
$dCell.Value = "'250.88"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("E24").Value = "  -2.36%  "
$dCell = $ws.Range("D25")
$dCell.Value = "'2.05"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$dCell = $ws.Range("D26")
$dCell.Value = "'26.42"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.02%  "
$dCell = $ws.Range("D29")
$dCell.Value = "'40.27"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +3.70%  "
$dCell = $ws.Range("D30")
$dCell.Value = "'10.37"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +2.09%  "
$dCell = $ws.Range("D31")
$dCell.Value = "'5.91"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$dCell = $ws.Range("D32")
$dCell.Value = "'156.42"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$dCell = $ws.Range("D33")
$dCell.Value = "'2.15"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E34").Value = "  +3.75%  "
$dCell = $ws.Range("D35")
$dCell.Value = "'3.32"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$dCell = $ws.Range("D36")
$dCell.Value = "'18.79"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -4.30%  "
$dCell = $ws.Range("D37")
$dCell.Value = "'0.0778"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  +7.74%  "
$dCell = $ws.Range("D41")
$dCell.Value = "'22.39"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -7.89%  "
$dCell = $ws.Range("D42")
$dCell.Value = "'3.81"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  -0.23%  "
$dCell = $ws.Range("D45")
$dCell.Value = "'2.032.89"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
$dCell = $ws.Range("D46")
$dCell.Value = "'3.24"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -4.63%  "
$dCell = $ws.Range("D47")
$dCell.Value = "'9.06"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +2.84%  "
$dCell = $ws.Range("D48")
$dCell.Value = "'84.43"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "
$dCell = $ws.Range("D49")
$dCell.Value = "'105.55"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "
$dCell = $ws.Range("D50")
$dCell.Value = "'75.04"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "
$dCell = $ws.Range("D51")
$dCell.Value = "'2.773.56"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
